$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the multiple-imputation standard errors (previously "(nan)" placeholders)
# for the theta_se (row 4) and lambda_se (row 6) rows, column by column.
$ws.Range("B4").Value = "(0.15)"
$ws.Range("B6").Value = "(0.17)"

$ws.Range("C4").Value = "(0.56)"
$ws.Range("C6").Value = "(0.16)"

$ws.Range("D4").Value = "(0.45)"
$ws.Range("D6").Value = "(0.11)"

$ws.Range("E4").Value = "(0.7)"
$ws.Range("E6").Value = "(0.35)"

$ws.Range("F4").Value = "(0.32)"
$ws.Range("F6").Value = "(0.25)"

$ws.Range("G4").Value = "(0.6)"
$ws.Range("G6").Value = "(0.4)"

$ws.Range("H4").Value = "(1.39)"
$ws.Range("I4").Value = "(1.39)"

$ws.Range("H6").Value = "(0.5)"
$ws.Range("I6").Value = "(1.45)"

$ws.Range("J4").Value = "(2.42)"
$ws.Range("J6").Value = "(0.88)"
